$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add an input/output header row (row 1 was blank) describing each logged
# column: input columns (Date/Time, Method, elapsedMs, wordCount,
# sentenceCount) and the analysis output columns (posWordCount .. the
# percentage/phrase text boxes).
$headers = @("Date/Time","Method","elapsedMs","wordCount","sentenceCount","posWordCount","negWordCount","posWordPercentage","negWordPercentage","posPhraseCount","negativePhraseCount","posWordPercentage","negPhrasePercentage")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# A1 picks up the same date/time display style already used by column A's
# data cells below it.
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Existing run (row 2) re-measured after switching the worker thread from a
# blocking Wait to Await -- timings and counts shift accordingly.
$ws.Cells.Item(2, 1).Value = 42585.690983796296
$ws.Cells.Item(2, 2).Value = "Bag"
$ws.Cells.Item(2, 3).Value = 8286
$ws.Cells.Item(2, 4).Value = 13068
$ws.Cells.Item(2, 5).Value = 1529
$ws.Cells.Item(2, 6).Value = 227
$ws.Cells.Item(2, 7).Value = 115
$ws.Cells.Item(2, 8).Value = 65
$ws.Cells.Item(2, 9).Value = 33
$ws.Cells.Item(2, 10).Value = 10
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 65
$ws.Cells.Item(2, 13).Value = 9

# Four additional Await-based runs logged right after.
$rows = @(
    @(42585.694803240738, "Bag", 8277, 13180, 1559, 231, 114, 64, 31, 3, 1, 64, 25),
    @(42585.698263888888, "Bag", 8157, 13177, 1559, 230, 114, 66, 33, 3, 1, 66, 25),
    @(42585.704293981478, "Bag", 8262, 13187, 1540, 230, 116, 64, 32, 10, 1, 64, 9),
    @(42585.70722222222, "Bag", 8494, 13313, 1571, 234, 115, 65, 32, 3, 1, 65, 25)
)

$r = 3
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# Give the new date/time cells (A3:A6) the same display style as A2.
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Resize every used column to fit its new (wider) header/content text.
$ws.Columns("A:M").AutoFit() | Out-Null
